# Apply the commit's data/formatting corrections to the workbook.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Costo"
$ws2 = $wb.Worksheets.Item(2)   # "Ultimo"

# -------------------------------------------------------------------
# 1) Sheet "Ultimo": raw-material / pallet unit prices were corrected
#    from absurd bulk totals down to realistic per-kg prices.
# -------------------------------------------------------------------
$ws2.Range("E1").Value  = 0.77    # A1224
$ws2.Range("E2").Value  = 0.62    # A0102
$ws2.Range("E3").Value  = 0.92    # A0404
$ws2.Range("E4").Value  = 0.325   # A0405
$ws2.Range("E5").Value  = 1.09    # A0601
$ws2.Range("E6").Value  = 4.4     # A2035
$ws2.Range("E8").Value  = 10.7    # B6107 pallet (LAH/004045)
$ws2.Range("E9").Value  = 0.645   # A1004
$ws2.Range("E10").Value = 0.68    # A1216
$ws2.Range("E11").Value = 0.98    # A0600
$ws2.Range("E12").Value = 2.4     # A3004
$ws2.Range("E13").Value = 0.88    # A0402
$ws2.Range("E14").Value = 0.92    # A0404
$ws2.Range("E15").Value = 0.31    # A0403
$ws2.Range("E17").Value = 10.7    # B6107 pallet (LAH/004047)

# -------------------------------------------------------------------
# 2) Sheet "Costo": recompute the dependent totals for both loads.
# -------------------------------------------------------------------
$ws1.Range("J2").Value = 0.98546
$ws1.Range("K2").Value = 1.006269135802469
$ws1.Range("L2").Value = -0.02080913580246913

$ws1.Range("J3").Value = 0.81349
$ws1.Range("K3").Value = 0.7661045950413224
$ws1.Range("L3").Value = 0.0473854049586776

# Long free-text detail cells recomputed with the corrected prices.
$text34 = @"
Lavorazioni toccate:
 [LAH/004045 q.: 4050.0]
Totale carichi: 4050.0

Lavorazione Linea 5: euro/kg. 0.264 x 4050.0 = 1069.2

Costi materie prime:
Lavoration LAH/004045:
 - A1224: EUR 0.77 x q. 2800.0 = 2156.0
 - A0102: EUR 0.62 x q. 240.0 = 148.8
 - A0404: EUR 0.92 x q. 420.0 = 386.4
 - A0405: EUR 0.325 x q. 420.0 = 136.5
 - A0601: EUR 1.09 x q. 32.0 = 34.88
 - A2035: EUR 4.4 x q. 8.0 = 35.2
Totale materie prime: 2897.78

Costi imballi e pallet:
 - Imballo [LAH/004045] B6003: EUR 0.405 x q. 162 = 65.61
 - Pallet [LAH/004045] B6107: EUR 10.7 x q. 4 = 42.8 
Totale imballi: 3006.19

Peso materie prime: 3920.0

Costo totale:
EUR 4075.39 : q. 4050.0 = EUR/unit 1.0062691358 (carico)

"@
$ws1.Range("I2").Value = $text34

$text48 = @"
Lavorazioni toccate:
 [LAH/004047 q.: 6050.0]
Totale carichi: 6050.0

Lavorazione Linea 4: euro/kg. 0.264 x 6050.0 = 1597.2

Costi materie prime:
Lavoration LAH/004047:
 - A1004: EUR 0.645 x q. 969.0 = 625.005
 - A1216: EUR 0.68 x q. 969.0 = 658.92
 - A0600: EUR 0.98 x q. 229.5 = 224.91
 - A3004: EUR 2.4 x q. 102.0 = 244.8
 - A0402: EUR 0.88 x q. 17.85 = 15.708
 - A0404: EUR 0.92 x q. 384.03 = 353.3076
 - A0403: EUR 0.31 x q. 2428.62 = 752.8722
Totale materie prime: 2875.5228

Costi imballi e pallet:
 - Imballo [LAH/004047] B6003: EUR 0.405 x q. 242 = 98.01
 - Pallet [LAH/004047] B6107: EUR 10.7 x q. 6 = 64.2 
Totale imballi: 3037.7328

Peso materie prime: 5100.0

Costo totale:
EUR 4634.9328 : q. 6050.0 = EUR/unit 0.766104595041 (carico)

"@
$ws1.Range("I3").Value = $text48

# New "Use history" marker for the second load, mirroring G3.
$ws1.Range("M3").Value = "X"

# Setting the long wrapped strings makes the engine auto-grow the row
# height; put it back to the sheet's standard height (no explicit
# per-row override), same as before the edit.
$ws1.Rows.Item(2).AutoFit()
$ws1.Rows.Item(3).AutoFit()

# -------------------------------------------------------------------
# 3) Formatting: the data rows on "Costo" now carry the same boxed /
#    centred style already used on "Ultimo" (font Arial 9, thin
#    border, centred horizontally) instead of being unformatted.
# -------------------------------------------------------------------
$dataRange = $ws1.Range("A2:M3")
$dataRange.Font.Name = "Arial"
$dataRange.Font.Size = 9
$dataRange.HorizontalAlignment = -4108   # xlCenter
$dataRange.Borders.LineStyle = 1         # xlContinuous

# -------------------------------------------------------------------
# 4) Both fonts used in the workbook (bold header + regular body)
#    were switched from Arial to Verdana.
# -------------------------------------------------------------------
$ws1.Range("A1:M1").Font.Name = "Verdana"
$ws1.Range("A2:M3").Font.Name = "Verdana"
$ws2.Range("A1:E17").Font.Name = "Verdana"

$ws1.Rows.Item(2).AutoFit()
$ws1.Rows.Item(3).AutoFit()

# -------------------------------------------------------------------
# 5) Explicit column widths were set on "Costo".
# -------------------------------------------------------------------
$ws1.Columns.Item(1).ColumnWidth  = 10
$ws1.Columns.Item(2).ColumnWidth  = 10
$ws1.Columns.Item(3).ColumnWidth  = 10
$ws1.Columns.Item(4).ColumnWidth  = 2
$ws1.Columns.Item(5).ColumnWidth  = 10
$ws1.Columns.Item(6).ColumnWidth  = 10
$ws1.Columns.Item(7).ColumnWidth  = 5
$ws1.Columns.Item(8).ColumnWidth  = 10
$ws1.Columns.Item(9).ColumnWidth  = 40
$ws1.Columns.Item(10).ColumnWidth = 15
$ws1.Columns.Item(11).ColumnWidth = 15
$ws1.Columns.Item(12).ColumnWidth = 15
$ws1.Columns.Item(13).ColumnWidth = 5
